$d = $word.ActiveDocument

# 1. "initial" -> "early"
$d.Content.Find.Execute("initial", $true, $false, $false, $false, $false, $true, 1, $false, "early", 2) | Out-Null

# 2. Insert "General Workflow" heading paragraph + numbered list after paragraph 4
#    (paragraph 4 ends "...waiting to smack something."; paragraph 5 is the blank spacer)
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter() | Out-Null

$pHeading = $d.Paragraphs(5)
$pHeading.Range.InsertBefore("General Workflow")
$pHeading2 = $d.Paragraphs(5)
$pHeading2.Range.InsertParagraphAfter() | Out-Null

# --- List item 1 ---
$p = $d.Paragraphs(6)
$p.Range.InsertBefore('From @URL1@, download @ZIP1@, extract all files and open "@EXE1@".')
$p.Range.ListFormat.ApplyNumberDefault()

$f = $p.Range.Duplicate
$f.Find.Execute("@URL1@") | Out-Null
$h = $d.Hyperlinks.Add($f, "https://github.com/rrskybox/TransientSearch/tree/master/publish", "", "", "https://github.com/rrskybox/TransientSearch/tree/master/publish")
$hr = $h.Range
$hr.Font.Name = "Arial"
$hr.Font.NameAscii = "Arial"
$hr.Font.Size = 10

$f = $p.Range.Duplicate
$f.Find.Execute("@ZIP1@") | Out-Null
$f.Text = "TransientSearch.zip"
$f.Italic = 1
$f.Font.ItalicBi = 1

$f = $p.Range.Duplicate
$f.Find.Execute("@EXE1@") | Out-Null
$f.Text = "setup.exe"
$f.Italic = 1
$f.Font.ItalicBi = 1

$p.Range.InsertParagraphAfter() | Out-Null

# --- List item 2 ---
$p = $d.Paragraphs(7)
$p.Range.InsertBefore('From https://github.com/rrskybox/Hot-Pursuit/tree/master/Hot%20Pursuit/publish, download @ZIP2@, extract all files and open "@EXE2@".')
$p.Range.ListFormat.ApplyNumberDefault()

$f = $p.Range.Duplicate
$f.Find.Execute("@ZIP2@") | Out-Null
$f.Text = "HotPursuit64.zip"
$f.Italic = 1
$f.Font.ItalicBi = 1

$f = $p.Range.Duplicate
$f.Find.Execute("@EXE2@") | Out-Null
$f.Text = "setup.exe"
$f.Italic = 1
$f.Font.ItalicBi = 1

$p.Range.InsertParagraphAfter() | Out-Null

# --- List item 3 ---
$p = $d.Paragraphs(8)
$p.Range.InsertBefore('Launch @SKY1@.  From the TSXToolKit in the Start Menu: Launch @TS1@.  Launch @HP1@.')
$p.Range.ListFormat.ApplyNumberDefault()

$f = $p.Range.Duplicate
$f.Find.Execute("@SKY1@") | Out-Null
$f.Text = "TheSky64"
$f.Italic = 1
$f.Font.ItalicBi = 1

$f = $p.Range.Duplicate
$f.Find.Execute("@TS1@") | Out-Null
$f.Text = "Transient Search"
$f.Italic = 1
$f.Font.ItalicBi = 1

$f = $p.Range.Duplicate
$f.Find.Execute("@HP1@") | Out-Null
$f.Text = "Hot Pursuit"
$f.Italic = 1
$f.Font.ItalicBi = 1

$p.Range.InsertParagraphAfter() | Out-Null

# --- List item 4 ---
$p = $d.Paragraphs(9)
$p.Range.InsertBefore('Transient Search: In the MPC NEO box, check "@SCOUT1@" and select "@NEO1@".')
$p.Range.ListFormat.ApplyNumberDefault()

$f = $p.Range.Duplicate
$f.Find.Execute("@SCOUT1@") | Out-Null
$f.Text = "Scout"
$f.Italic = 1
$f.Font.ItalicBi = 1

$f = $p.Range.Duplicate
$f.Find.Execute("@NEO1@") | Out-Null
$f.Text = "NEO"
$f.Italic = 1
$f.Font.ItalicBi = 1

$p.Range.InsertParagraphAfter() | Out-Null

# --- List item 5 ---
$p = $d.Paragraphs(10)
$p.Range.InsertBefore('TheSky64: @EDIT1@.  Pick a NEO asteroid from the Sky Chart.')
$p.Range.ListFormat.ApplyNumberDefault()

$f = $p.Range.Duplicate
$f.Find.Execute("@EDIT1@") | Out-Null
$f.Text = "Edit->Paste Photo"
$f.Italic = 1
$f.Font.ItalicBi = 1

$p.Range.InsertParagraphAfter() | Out-Null

# --- List item 6 ---
$p = $d.Paragraphs(11)
$p.Range.InsertBefore('Hot Pursuit:  "@PURSUE1@".')
$p.Range.ListFormat.ApplyNumberDefault()

$f = $p.Range.Duplicate
$f.Find.Execute("@PURSUE1@") | Out-Null
$f.Text = "Pursue"
$f.Italic = 1
$f.Font.ItalicBi = 1

$p.Range.InsertParagraphAfter() | Out-Null

# --- List item 7 ---
$p = $d.Paragraphs(12)
$p.Range.InsertBefore('TheSky64: Image to taste.')
$p.Range.ListFormat.ApplyNumberDefault()

Write-Output "done with workflow section"
